$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.123.60'
$ws.Range("E2").Value = '  +4.61%  '
$ws.Range("D3").Value = '2.247.74'
$ws.Range("E3").Value = '  +3.80%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''244.99'
$ws.Range("E5").Value = '  +3.79%  '
$ws.Range("D6").Value = '''0.616'
$ws.Range("E6").Value = '  +1.67%  '
$ws.Range("D7").Value = '''75.30'
$ws.Range("E7").Value = '  +8.54%  '
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("D9").Value = '''0.605'
$ws.Range("E9").Value = '  +6.78%  '
$ws.Range("D10").Value = '''41.31'
$ws.Range("D11").Value = '''0.0932'
$ws.Range("E11").Value = '  +2.69%  '
$ws.Range("D12").Value = '''6.95'
$ws.Range("E12").Value = '  +4.21%  '
$ws.Range("D13").Value = '''0.102'
$ws.Range("E13").Value = '  +1.59%  '
$ws.Range("D14").Value = '2.585.79'
$ws.Range("E14").Value = '  +3.88%  '
$ws.Range("D15").Value = '''14.60'
$ws.Range("E15").Value = '  +1.96%  '
$ws.Range("D16").Value = '2.241.95'
$ws.Range("E16").Value = '  +4.09%  '
$ws.Range("D17").Value = '''0.793'
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("D18").Value = '43.029.08'
$ws.Range("E18").Value = '  +4.82%  '
$ws.Range("E19").Value = '  +5.77%  '
$ws.Range("D20").Value = '''71.13'
$ws.Range("E20").Value = '  +1.94%  '
$ws.Range("D21").Value = '''5.99'
$ws.Range("E21").Value = '  +3.49%  '
$ws.Range("D22").Value = '''9.89'
$ws.Range("E22").Value = '  +6.29%  '
$ws.Range("D23").Value = '''229.77'
$ws.Range("E23").Value = '  +2.27%  '
$ws.Range("E24").Value = '  +16.97%  '
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("D26").Value = '''10.88'
$ws.Range("E26").Value = '  +2.32%  '
$ws.Range("D27").Value = '''3.41'
$ws.Range("E27").Value = '  +3.91%  '
$ws.Range("E28").Value = '  +2.83%  '
$ws.Range("D29").Value = '''38.42'
$ws.Range("E29").Value = '  +27.61%  '
$ws.Range("E30").Value = '  +2.42%  '
$ws.Range("D31").Value = '''171.62'
$ws.Range("E31").Value = '  +1.90%  '
$ws.Range("E32").Value = '  +2.57%  '
$ws.Range("D33").Value = '''0.0797'
$ws.Range("E33").Value = '  +5.79%  '
$ws.Range("D34").Value = '''5.30'
$ws.Range("E34").Value = '  +4.30%  '
$ws.Range("E35").Value = '  +1.83%  '
$ws.Range("E36").Value = '  +7.15%  '
$ws.Range("D37").Value = '''4.37'
$ws.Range("E37").Value = '  +6.96%  '
$ws.Range("E38").Value = '  +18.98%  '
$ws.Range("D39").Value = '''13.09'
$ws.Range("E39").Value = '  +14.33%  '
$ws.Range("D40").Value = '''2.13'
$ws.Range("E40").Value = '  +3.60%  '
$ws.Range("D41").Value = '''0.205'
$ws.Range("E41").Value = '  +10.00%  '
$ws.Range("D42").Value = '''5.44'
$ws.Range("E42").Value = '  +2.67%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '''105.33'
$ws.Range("E43").Value = '  +9.00%  '
$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").Value = '''59.41'
$ws.Range("E44").Value = '  +2.40%  '
$ws.Range("D45").Value = '''8.71'
$ws.Range("E45").Value = '  +5.90%  '
$ws.Range("D46").Value = '''0.483'
$ws.Range("E46").Value = '  +30.93%  '
$ws.Range("E47").Value = '  +3.17%  '
$ws.Range("D48").Value = '''2.40'
$ws.Range("E48").Value = '  +11.36%  '
$ws.Range("E49").Value = '  +3.39%  '
$ws.Range("E50").Value = '  +3.38%  '
$ws.Range("D51").Value = '2.458.76'
$ws.Range("E51").Value = '  +3.97%  '
